$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# --- Weekly statistics table updates ---
$ws.Range("M14").Value = -75
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 66
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = -16.455696202531
$ws.Range("L16").Value = 127.586206896552
$ws.Range("M16").Value = -2.941176470588
$ws.Range("N16").Value = -76
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -28.571428571428
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = -22.222222222222
$ws.Range("I17").Value = 95
$ws.Range("J17").Value = 90
$ws.Range("K17").Value = 5.555555555555
$ws.Range("L17").Value = 46.153846153846
$ws.Range("M17").Value = 79.245283018867
$ws.Range("N17").Value = -19.491525423728
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 8.333333333333
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = 17.391304347826
$ws.Range("L18").Value = 125
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = -84.210526315789
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 71.428571428571
$ws.Range("F19").Value = 62
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = 31.914893617021
$ws.Range("I19").Value = 300
$ws.Range("J19").Value = 279
$ws.Range("K19").Value = 7.526881720430
$ws.Range("L19").Value = 63.043478260869
$ws.Range("M19").Value = 97.368421052631
$ws.Range("N19").Value = 54.639175257732
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 28.571428571428
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 19.047619047619
$ws.Range("I20").Value = 142
$ws.Range("J20").Value = 128
$ws.Range("K20").Value = 10.9375
$ws.Range("L20").Value = 222.727272727273
$ws.Range("M20").Value = 136.666666666667
$ws.Range("N20").Value = -82.447466007416
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 13.793103448275
$ws.Range("F21").Value = 134
$ws.Range("G21").Value = 125
$ws.Range("H21").Value = 7.2
$ws.Range("I21").Value = 691
$ws.Range("J21").Value = 657
$ws.Range("K21").Value = 5.175038051750
$ws.Range("L21").Value = 89.315068493150
$ws.Range("M21").Value = 74.936708860759
$ws.Range("N21").Value = -64.010416666666
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("G22").NumberFormat = '#,##0'
$ws.Range("H22").Value = 0
$ws.Range("H22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = 50
$ws.Range("M22").Value = -40
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 200
$ws.Range("F23").Value = 11
$ws.Range("H23").Value = 22.222222222222
$ws.Range("I23").Value = 39
$ws.Range("J23").Value = 37
$ws.Range("K23").Value = 5.405405405405
$ws.Range("L23").Value = 116.666666666667
$ws.Range("M23").Value = 85.714285714285
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = -9.302325581395
$ws.Range("I24").Value = 473
$ws.Range("J24").Value = 485
$ws.Range("K24").Value = -2.474226804123
$ws.Range("L24").Value = 29.234972677595
$ws.Range("M24").Value = 30.303030303030
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -61.538461538461
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -31.428571428571
$ws.Range("I25").Value = 141
$ws.Range("J25").Value = 151
$ws.Range("K25").Value = -6.622516556291
$ws.Range("L25").Value = 31.775700934579
$ws.Range("M25").Value = -19.428571428571
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = -16.666666666666
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("L28").Value = -66.666666666666
$ws.Range("M28").Value = -33.333333333333
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("L29").Value = -60
$ws.Range("M29").Value = -33.333333333333
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = '#,##0'
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = '#,##0'
$ws.Range("H30").Value = 0
$ws.Range("H30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -66.666666666666
